# Apply the cryptos list update (prices + 1h volume % changes, plus
# the EnergySwap/Maker and VeChain/OKB row swaps) for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.762.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.184.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.735.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.798.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.206.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.19%  "

$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0872"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("E33").Value = "  +2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.18%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.79%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.785.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.47%  "

$ws.Range("E39").Value = "  +2.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0295"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.95%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.713"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.226.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.984"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.798"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("E51").Value = "  +0.01%  "
